$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunnerConfig")

# Update existing rows 2-4: ToRun column changes from Y to N
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"

# Add two new test case rows (6 and 7)
$ws.Range("B6").Value = "TestCase 01"
$ws.Range("B7").Value = "Test Case02 "
$ws.Range("A6").Value = "testCase01"
$ws.Range("A7").Value = "testCase02"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Match the bordered formatting used by the row above (B5:C5)
$ws.Range("B5").Copy()
$ws.Range("A6:C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to reflect the new active cell
$ws.Range("B11").Select()
